# Auto-applied numeric corrections to Marilith_Profits workbook (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 255  # was 303.33334
$ws.Range("I11").Value = 255  # was 303.33334
$ws.Range("K11").Value = 255  # was 303.33334
$ws.Range("M11").Value = -115  # was -163.33334
$ws.Range("H32").Value = 833.25  # was 1028.4286
$ws.Range("I32").Value = 583.1667  # was 749.5
$ws.Range("J32").Value = 1083.3334  # was 1140
$ws.Range("K32").Value = 583.1667  # was 749.5
$ws.Range("L32").Value = 1083.3334  # was 1140
$ws.Range("M32").Value = -257.1667  # was -423.5
$ws.Range("N32").Value = -1735.3334  # was -1792
$ws.Range("H33").Value = 429.85715  # was 399
$ws.Range("I33").Value = 334.83334  # was 313.14285
$ws.Range("K33").Value = 334.83334  # was 313.14285
$ws.Range("M33").Value = -105.83334  # was -84.14285000000001
$ws.Range("H41").Value = 1539.1111  # was 1595.5454
$ws.Range("J41").Value = 2081.6667  # was 1988.8
$ws.Range("L41").Value = 2081.6667  # was 1988.8
$ws.Range("N41").Value = -2961.6667  # was -2868.8
$ws.Range("H64").Value = 5000  # was 4500.3335
$ws.Range("I64").Value = 3500  # was 2750.5
$ws.Range("J64").Value = 5750  # was 8000
$ws.Range("K64").Value = 3500  # was 2750.5
$ws.Range("L64").Value = 5750  # was 8000
$ws.Range("M64").Value = -3252  # was -2502.5
$ws.Range("N64").Value = -6246  # was -8496
$ws.Range("H67").Value = 5000  # was 4500.3335
$ws.Range("I67").Value = 3500  # was 2750.5
$ws.Range("J67").Value = 5750  # was 8000
$ws.Range("K67").Value = 3500  # was 2750.5
$ws.Range("L67").Value = 5750  # was 8000
$ws.Range("M67").Value = -2642  # was -1892.5
$ws.Range("N67").Value = -7466  # was -9716
$ws.Range("H104").Value = 228.42857  # was 228.57143
$ws.Range("I104").Value = 228.42857  # was 228.57143
$ws.Range("K104").Value = 685.28571  # was 685.71429
$ws.Range("M104").Value = 1061.71429  # was 1061.28571
$ws.Range("H107").Value = 1196.2593  # was 1209.1852
$ws.Range("I107").Value = 870.8333  # was 893.5217
$ws.Range("J107").Value = 3799.6667  # was 3024.25
$ws.Range("K107").Value = 870.8333  # was 893.5217
$ws.Range("L107").Value = 3799.6667  # was 3024.25
$ws.Range("M107").Value = 1049.1667  # was 1026.4783
$ws.Range("N107").Value = -7639.6667  # was -6864.25
$ws.Range("H116").Value = 25666  # was 12392.733
$ws.Range("I116").Value = 7997  # was 3821.3333
$ws.Range("J116").Value = 29199.8  # was 25249.834
$ws.Range("K116").Value = 7997  # was 3821.3333
$ws.Range("L116").Value = 29199.8  # was 25249.834
$ws.Range("M116").Value = -4555  # was -379.3332999999998
$ws.Range("N116").Value = -36083.8  # was -32133.834
$ws.Range("H131").Value = 3543.6924  # was 3028.875
$ws.Range("I131").Value = 3322.3333  # was 2961.7856
$ws.Range("J131").Value = 6200  # was 3498.5
$ws.Range("K131").Value = 9966.999899999999  # was 8885.356800000001
$ws.Range("L131").Value = 18600  # was 10495.5
$ws.Range("M131").Value = -4926.999899999999  # was -3845.356800000001
$ws.Range("N131").Value = -28680  # was -20575.5
$ws.Range("H137").Value = 3275  # was 2580.6155
$ws.Range("I137").Value = 2883.3333  # was 2068.625
$ws.Range("J137").Value = 3666.6667  # was 3399.8
$ws.Range("K137").Value = 8649.999899999999  # was 6205.875
$ws.Range("L137").Value = 11000.0001  # was 10199.4
$ws.Range("M137").Value = -6099.999899999999  # was -3655.875
$ws.Range("N137").Value = -16100.0001  # was -15299.4
$ws.Range("H138").Value = 4546.095  # was 4411.2
$ws.Range("J138").Value = 5076.4707  # was 5181.8184
$ws.Range("L138").Value = 15229.4121  # was 15545.4552
$ws.Range("N138").Value = -25509.4121  # was -25825.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24813.637  # was 22040.682
$ws.Range("I32").Value = 22185.715  # was 22040.682
$ws.Range("J32").Value = 80000  # was 0
$ws.Range("K32").Value = 22185.715  # was 22040.682
$ws.Range("L32").Value = 80000  # was 0
$ws.Range("M32").Value = -21898.715  # was -21753.682
$ws.Range("N32").Value = -80574  # was None
$ws.Range("H45").Value = 3000  # was 2144.5715
$ws.Range("I45").Value = 0  # was 2002
$ws.Range("K45").Value = 0  # was 2002
$ws.Range("M45").ClearContents()  # was -1625
$ws.Range("H102").Value = 2222.4285  # was 2285.6924
$ws.Range("I102").Value = 2222.4285  # was 2285.6924
$ws.Range("K102").Value = 2222.4285  # was 2285.6924
$ws.Range("M102").Value = -600.4285  # was -663.6923999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14998.75  # was 14997
$ws.Range("J35").Value = 14998.75  # was 14997
$ws.Range("L35").Value = 14998.75  # was 14997
$ws.Range("N35").Value = -15618.75  # was -15617

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5435.5  # was 5144.4287
$ws.Range("J31").Value = 6271.6665  # was 5553.25
$ws.Range("L31").Value = 6271.6665  # was 5553.25
$ws.Range("N31").Value = -6861.6665  # was -6143.25
$ws.Range("H34").Value = 5435.5  # was 5144.4287
$ws.Range("J34").Value = 6271.6665  # was 5553.25
$ws.Range("L34").Value = 6271.6665  # was 5553.25
$ws.Range("N34").Value = -6675.6665  # was -5957.25
$ws.Range("H59").Value = 30278.354  # was 29983.188
$ws.Range("J59").Value = 34997.5  # was 34997.11
$ws.Range("L59").Value = 34997.5  # was 34997.11
$ws.Range("N59").Value = -37287.5  # was -37287.11
$ws.Range("H102").Value = 0  # was 21900
$ws.Range("I102").Value = 0  # was 21900
$ws.Range("K102").Value = 0  # was 21900
$ws.Range("M102").ClearContents()  # was -19466
$ws.Range("H107").Value = 644.93335  # was 658.8570999999999
$ws.Range("I107").Value = 582.8889  # was 599.5
$ws.Range("K107").Value = 582.8889  # was 599.5
$ws.Range("M107").Value = 1337.1111  # was 1320.5
$ws.Range("H122").Value = 3181.7646  # was 3240.0625
$ws.Range("I122").Value = 3068.125  # was 3122.7334
$ws.Range("K122").Value = 9204.375  # was 9368.200199999999
$ws.Range("M122").Value = -6754.375  # was -6918.200199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 637.1667  # was 612.3461
$ws.Range("I5").Value = 515.65  # was 497.36365
$ws.Range("K5").Value = 1546.95  # was 1492.09095
$ws.Range("M5").Value = -1434.95  # was -1380.09095
$ws.Range("H21").Value = 416.66666  # was 124
$ws.Range("I21").Value = 225  # was 400
$ws.Range("J21").Value = 800  # was 55
$ws.Range("K21").Value = 675  # was 1200
$ws.Range("L21").Value = 2400  # was 165
$ws.Range("M21").Value = -502  # was -1027
$ws.Range("N21").Value = -2746  # was -511
$ws.Range("H135").Value = 637.1667  # was 612.3461
$ws.Range("I135").Value = 515.65  # was 497.36365
$ws.Range("K135").Value = 4640.849999999999  # was 4476.27285
$ws.Range("M135").Value = -2105.849999999999  # was -1941.27285
$ws.Range("H137").Value = 25676.666  # was 269257.5
$ws.Range("I137").Value = 40030  # was 520015
$ws.Range("K137").Value = 120090  # was 1560045
$ws.Range("M137").Value = -114990  # was -1554945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0  # was 500
$ws.Range("J25").Value = 0  # was 500
$ws.Range("L25").Value = 0  # was 500
$ws.Range("N25").ClearContents()  # was -1558
$ws.Range("H46").Value = 16000  # was 13600
$ws.Range("H132").Value = 1274  # was 1047.6
$ws.Range("I132").Value = 1274  # was 1047.6
$ws.Range("K132").Value = 3822  # was 3142.8
$ws.Range("M132").Value = -1292  # was -612.7999999999997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1698.7  # was 1810.1111
$ws.Range("I7").Value = 1698.7  # was 1810.1111
$ws.Range("K7").Value = 1698.7  # was 1810.1111
$ws.Range("M7").Value = -1586.7  # was -1698.1111
$ws.Range("H126").Value = 1698.7  # was 1810.1111
$ws.Range("I126").Value = 1698.7  # was 1810.1111
$ws.Range("K126").Value = 5096.1  # was 5430.3333
$ws.Range("M126").Value = -2626.1  # was -2960.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0  # was 6000
$ws.Range("J62").Value = 0  # was 6000
$ws.Range("L62").Value = 0  # was 6000
$ws.Range("N62").ClearContents()  # was -7248
$ws.Range("H65").Value = 0  # was 6000
$ws.Range("J65").Value = 0  # was 6000
$ws.Range("L65").Value = 0  # was 30000
$ws.Range("N65").ClearContents()  # was -36240
$ws.Range("H111").Value = 48721.75  # was 51629
$ws.Range("J111").Value = 48721.75  # was 51629
$ws.Range("L111").Value = 48721.75  # was 51629
$ws.Range("N111").Value = -56901.75  # was -59809
